$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Sheet1 "Air Movement " ---
# Row 1 headers: drop the old "Factor" (A1) / "Weight" (B1) labels, shift the
# pump-type headers one column to the left, and refresh the Air Pump header.
$ws.Range("A1").ClearContents()
$ws.Range("B1").Value = "Microblower"
$ws.Range("C1").Value = "Microfan"
$ws.Range("D1").Value = "Air Pump " + [char]10 + "(CurieJet)"
$ws.Range("D1").WrapText = $true
$ws.Range("E1").Value = "Passive"
$ws.Range("F1").Value = "Piezoelectric Pump"
$ws.Range("G1").Value = "Electromagnetic Pump"
$ws.Range("H1").ClearContents()

# Row labels
$ws.Range("A2").Value = "Power Consumption (W)"
$ws.Range("A3").Value = "Size"
$ws.Range("A4").Value = "Noise"
$ws.Range("A5").Value = "Durability"
$ws.Range("A6").Value = "Maintenance"
$ws.Range("A7").Value = "Airflow Rate"
$ws.Range("A8").Value = "Weight"
$ws.Range("A9").Value = "Link"

# Hyperlink cell
$ws.Range("D9").Value = "https://www.curiejet.com/en/product/micro-pump/air-pump-and-micro-blower"
$ws.Hyperlinks.Add($ws.Range("D9"), "https://www.curiejet.com/en/product/micro-pump/air-pump-and-micro-blower")

# Column widths (matches post-edit autofit state)
$ws.Columns.Item(1).ColumnWidth = 22.42578125
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(3).ColumnWidth = 8.85546875
$ws.Columns.Item(6).ColumnWidth = 18
$ws.Columns.Item(7).ColumnWidth = 21.140625

$ws.Range("E1").Select()

Write-Host "sheet1 done"
